# Applies the "Updated symbol list" GitHub Actions commit: refreshed
# price/volume figures in columns D/E, plus a re-sort of the exchange-
# token block (rows 18-25) that shifted B/C/D/E down by one row and
# inserted HotbitToken at the top of that block.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value. D/E columns hold numeric-looking
# text ("328.03", "-0.96%"), so NumberFormat is forced to Text ("@")
# before the write (otherwise Excel would coerce them to numbers), and
# the style is reset to Normal afterwards so no stray formatting sticks.
$changes = @(
    @{Cell='D2'; Value='328.03'}
    @{Cell='E2'; Value='-0.96%'}
    @{Cell='D3'; Value='43.91'}
    @{Cell='E3'; Value='5.57%'}
    @{Cell='D4'; Value='5.408'}
    @{Cell='E4'; Value='-5.00%'}
    @{Cell='D5'; Value='0.08094'}
    @{Cell='E5'; Value='-3.09%'}
    @{Cell='D6'; Value='8.673'}
    @{Cell='E6'; Value='-1.51%'}
    @{Cell='D7'; Value='1.901'}
    @{Cell='E7'; Value='-5.46%'}
    @{Cell='D8'; Value='4.299'}
    @{Cell='E8'; Value='-3.72%'}
    @{Cell='E9'; Value='-5.06%'}
    @{Cell='D10'; Value='0.9419'}
    @{Cell='E10'; Value='1.77%'}
    @{Cell='D11'; Value='0.1185'}
    @{Cell='E11'; Value='-7.98%'}
    @{Cell='E12'; Value='-3.88%'}
    @{Cell='D13'; Value='0.09597'}
    @{Cell='E13'; Value='0.77%'}
    @{Cell='D14'; Value='0.04219'}
    @{Cell='E14'; Value='9.67%'}
    @{Cell='D15'; Value='0.1069'}
    @{Cell='E15'; Value='0.71%'}
    @{Cell='D16'; Value='0.001288'}
    @{Cell='E16'; Value='-1.20%'}
    @{Cell='D17'; Value='0.005983'}
    @{Cell='E17'; Value='-2.07%'}
    @{Cell='B18'; Value='HotbitToken'}
    @{Cell='C18'; Value='https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'}
    @{Cell='D18'; Value='0.004308'}
    @{Cell='E18'; Value='-1.42%'}
    @{Cell='B19'; Value='LEO'}
    @{Cell='C19'; Value='https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'}
    @{Cell='D19'; Value='3.553'}
    @{Cell='E19'; Value='3.42%'}
    @{Cell='B20'; Value='BitpandaEcosystemToken'}
    @{Cell='C20'; Value='https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'}
    @{Cell='D20'; Value='0.3517'}
    @{Cell='E20'; Value='-0.51%'}
    @{Cell='B21'; Value='MCDex'}
    @{Cell='C21'; Value='https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'}
    @{Cell='D21'; Value='8.508'}
    @{Cell='E21'; Value='-1.09%'}
    @{Cell='B22'; Value='ProBitToken'}
    @{Cell='C22'; Value='https://coinranking.com/coin/lQP4d6T2+probittoken-prob'}
    @{Cell='D22'; Value='0.1359'}
    @{Cell='E22'; Value='-0.27%'}
    @{Cell='B23'; Value='ZBToken'}
    @{Cell='C23'; Value='https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'}
    @{Cell='D23'; Value='0.2608'}
    @{Cell='E23'; Value='3.89%'}
    @{Cell='B24'; Value='CoinExToken'}
    @{Cell='C24'; Value='https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'}
    @{Cell='D24'; Value='0.04368'}
    @{Cell='E24'; Value='-0.72%'}
    @{Cell='B25'; Value='BitKan'}
    @{Cell='C25'; Value='https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'}
    @{Cell='D25'; Value='0.001240'}
    @{Cell='E25'; Value='-2.65%'}
    @{Cell='E26'; Value='1.38%'}
    @{Cell='D27'; Value='0.0004015'}
    @{Cell='E27'; Value='0.56%'}
    @{Cell='D39'; Value='0.02700'}
    @{Cell='E39'; Value='-4.13%'}
    @{Cell='D40'; Value='0.05497'}
    @{Cell='E40'; Value='-0.52%'}
    @{Cell='D41'; Value='0.007799'}
    @{Cell='E41'; Value='-1.94%'}
    @{Cell='D42'; Value='0.009759'}
    @{Cell='E42'; Value='5.69%'}
    @{Cell='D43'; Value='0.1395'}
    @{Cell='E43'; Value='-2.62%'}
    @{Cell='D44'; Value='0.002140'}
    @{Cell='E44'; Value='3.78%'}
    @{Cell='D45'; Value='0.009634'}
    @{Cell='E45'; Value='-17.93%'}
    @{Cell='D46'; Value='0.00007107'}
    @{Cell='E46'; Value='3.11%'}
    @{Cell='D47'; Value='0.00000000755'}
    @{Cell='E47'; Value='0.56%'}
    @{Cell='D48'; Value='0.003474'}
    @{Cell='E48'; Value='0.31%'}
    @{Cell='D49'; Value='0.002285'}
    @{Cell='E49'; Value='0.22%'}
    @{Cell='D50'; Value='0.00002114'}
    @{Cell='E50'; Value='0.56%'}
    @{Cell='D51'; Value='0.0002013'}
    @{Cell='E51'; Value='0.56%'}
)

foreach ($change in $changes) {
    $range = $ws.Range($change.Cell)
    if ($change.Cell -match "^[DE]\d+$") {
        $range.NumberFormat = "@"
        $range.Value = $change.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $change.Value
    }
}
